$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (IA Control) list reordering / de-duplication fixes ---
$ws.Range("A2").Value = "AU-4,AU-4 (1)"
$ws.Range("A3").Value = "SC-5 (2),CM-6 b,SC-5"
$ws.Range("A4").Value = "AU-7 a,CM-5 (1),AU-7 b,AU-12 (3),AC-6 (9),AU-8 b,AC-6 (8)"
$ws.Range("A5").Value = "AC-17 (9),CM-6 b,CM-7 b,AC-17 (1)"
$ws.Range("A8").Value = "IA-2 (11),IA-2 (12)"
$ws.Range("A10").Value = "CM-7 (2),CM-7 (5) (b)"
$ws.Range("A12").Value = "AC-7 b,AC-7 a"
$ws.Range("A15").Value = "AU-3 (1),IA-2,IA-8"
$ws.Range("A16").Value = "CM-6 b,AC-6 (10)"
$ws.Range("A17").Value = "MA-4 (1) (a),AU-12 c,AU-12 a,AU-3,AU-3 (1)"
$ws.Range("A19").Value = "IA-5 (1) (b),IA-5 (1) (a),CM-6 b"
$ws.Range("A21").Value = "AC-12,SC-10,MA-4 (7),MA-4 e"
$ws.Range("A22").Value = "MA-4 (1) (a),AU-7 a,CM-5 (1),AU-12 a,AU-7 (1),AU-6 (4),AU-3,CM-6 b,AU-14 (1),AU-3 (1)"
$ws.Range("A25").Value = "MA-4 (1) (a),AU-12 c,AU-12 a,AU-3,AU-3 (1)"
$ws.Range("A29").Value = "SC-8 (2),SC-8,SC-8 (1)"
$ws.Range("A31").Value = "MA-4 (1) (a),AU-12 c,AU-12 a,AC-2 (4),AU-3,AU-3 (1)"
$ws.Range("A34").Value = "AC-11 b,AC-11 a"
$ws.Range("A38").Value = "SI-11 b,AU-9"
$ws.Range("A39").Value = "AU-3,CM-6 b"
$ws.Range("A41").Value = "AU-4 (1)"
$ws.Range("A42").Value = "SC-28 (1),SC-28"
$ws.Range("A44").Value = "IA-11"
$ws.Range("A45").Value = "AC-8 a,AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 b"
$ws.Range("A48").Value = "CM-6 b"
$ws.Range("A50").Value = "IA-2 (5),CM-6 b"
$ws.Range("A53").Value = "SC-13,MA-4 (6)"
$ws.Range("A55").Value = "AC-17 (2),SC-8"
$ws.Range("A56").Value = "MA-4 (1) (a),AU-12 c"
$ws.Range("A63").Value = "AU-5 a,AU-5 (1)"
$ws.Range("A65").Value = "IA-2 (2),CM-6 b"
$ws.Range("A67").Value = "MA-4 (1) (a),AU-12 c,AU-12 a,AU-3,AU-3 (1)"
$ws.Range("A69").Value = "CM-5 (1),AU-7 a,AU-12 c,AU-12 (3),AU-7 b,AU-12 a,AU-8 b,CM-6 b"
$ws.Range("A71").Value = "AU-3,AU-4 (1)"
$ws.Range("A77").Value = "MA-4 (1) (a),AU-12 c,AU-12 a,AC-2 (4),AU-3,AU-3 (1)"
$ws.Range("A79").Value = "AU-9 (3),AU-9"
$ws.Range("A80").Value = "IA-2 (1),IA-2 (2),IA-2 (4),IA-2 (3)"
$ws.Range("A81").Value = "CM-5 (3),CM-6 b"
$ws.Range("A86").Value = "MA-4 (1) (a),AU-12 c,AU-12 a,AU-3,AU-3 (1)"
$ws.Range("A88").Value = "AC-2 (4),AC-6 (9),CM-5 (1),AU-12 c"
$ws.Range("A89").Value = "IA-2 (5),IA-2 (3),IA-2 (2),IA-2 (4),IA-2"
$ws.Range("A90").Value = "IA-2 (11),IA-2 (12)"
$ws.Range("A91").Value = "AU-9 (3),AU-9"
$ws.Range("A96").Value = "AC-18 (1),SC-8,SC-8 (1)"
$ws.Range("A97").Value = "AU-8 (1) (b),AU-8 (1) (a),AU-8 b"
$ws.Range("A99").Value = "AU-9"
$ws.Range("A101").Value = "AC-3 (4),IA-11"
$ws.Range("A102").Value = "MA-4 (1) (a),AU-12 c,AU-12 a,AU-3,AU-3 (1)"
$ws.Range("A111").Value = "AU-5 a,AU-5 b"
$ws.Range("A119").Value = "MA-4 (1) (a),AU-12 c,AU-12 a,AU-3,AU-3 (1)"
$ws.Range("A123").Value = "CM-7 a,CM-7 b"
$ws.Range("A124").Value = "MA-4 (1) (a),AU-12 c,AU-12 a,AU-3,AU-3 (1)"
$ws.Range("A125").Value = "AC-18 (1),CM-7 a"
$ws.Range("A128").Value = "CM-7 a,CM-6 b,IA-5 (1) (c)"
$ws.Range("A136").Value = "AC-11 (1),AC-11 b"
$ws.Range("A139").Value = "CM-3 (5),SI-6 d,SI-6 b"
$ws.Range("A142").Value = "AC-2 (2)"
$ws.Range("A148").Value = "MA-4 (1) (a),AU-12 c,AU-12 a,AU-3,AU-14 (1),AU-3 (1)"
$ws.Range("A157").Value = "MA-4 (1) (a),AU-12 c,AU-12 a,AU-3,AU-3 (1)"
$ws.Range("A159").Value = "AC-17 (2),SC-8"
$ws.Range("A175").Value = "CM-7 a,SI-16"
$ws.Range("A181").Value = "SC-3,CM-6 b"

# --- Other per-cell content updates ---
# Row 42
$ws.Range("M42").Value = "Configure Red Hat Enterprise Linux 9 to prevent unauthorized modification of all information at rest by using disk encryption.`nEncrypting a partition in an already installed system is more difficult, because existing partitions will need to be resized and changed.`nTo encrypt an entire partition, dedicate a partition for encryption in the partition layout."
# Row 100
$ws.Range("F100").Value = "Red Hat Enterprise Linux 9 must allow the use of a temporary password for system logons with an immediate change to a permanent password."
$ws.Range("H100").Value = "Without providing this capability, an account may be created without a password.`nNon-repudiation cannot be guaranteed once an account is created if a user is not forced to change the temporary password upon initial logon.`nTemporary passwords are typically used to allow access when new accounts are created or passwords are changed.`nIt is common practice for administrators to create temporary passwords for user accounts that allow the users to log on, yet force them to change the password once they have successfully authenticated."
$ws.Range("I100").Value = "Applicable - Inherently Met"
$ws.Range("K100").Value = "Red Hat Enterprise Linux 9 supports this requirement and cannot be configured to be out of compliance.`nRed Hat Enterprise Linux 9 inherently meets this requirement."
$ws.Range("P100").Value = "Red Hat Enterprise Linux 9 offers the following commands to facilitate the use of a temporary password.`nchage -d 0 [username]`n(forces the user to change their password at next logon)`npasswd -e [username]`n(expires the passwd for a given user forcing a change at next logon.)"
$ws.Range("Q100").Value = "Red Hat Enterprise Linux 9 has the capability to perform temporary passwords based on organization policy.`nConfiguration is not appropriate to define at an enterprise level."
# Row 103
$ws.Range("F103").Value = "Red Hat Enterprise Linux 9 system must implement replay-resistant authentication mechanisms for network access to non-privileged accounts."
$ws.Range("I103").Value = "Applicable - Inherently Met"
$ws.Range("K103").Value = "Red Hat Enterprise Linux 9 supports this requirement and cannot be configured to be out of compliance.`nRed Hat Enterprise Linux 9 inherently meets this requirement."
$ws.Range("M103").Value = "Red Hat Enterprise Linux 9 inherently meets this requirement.`nNo fix is required."
$ws.Range("P103").Value = "The release notes of OpenSSH 7.6 states `"OpenSSH is a 100% complete SSH protocol 2.0 implementation and includes sftp client and server support.`"`nhttps://www.openssh.com/txt/release-7.6"
$ws.Range("Q103").Value = "The OpenSSH package in Red Hat Enterprise Linux 9 is version 8.7, which is newer than 7.6 which only supports SSH protocol 2.0 which is restraint to replay attacks."
# Row 192
$ws.Range("F192").Value = "Red Hat Enterprise Linux 9 must protect the confidentiality and integrity of all information at rest."
